$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A18").Value = "Umami"
$ws.Range("A18").Select()
